# Update TPM-derived statistics for the Cadm1-Cadm3 LR-pairs sheet
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.019613333333333
$ws.Range("H2").Value = 12.05884
$ws.Range("I2").Value = 0.4377217086785624
$ws.Range("J2").Value = 0.4377217086785624
$ws.Range("M2").Value = 1.328826
$ws.Range("N2").Value = 3.986478
$ws.Range("O2").Value = 0.03246115949735631
$ws.Range("P2").Value = 0.03246115949735631
$ws.Range("Q2").Value = 5.34136670728
$ws.Range("R2").Value = 48.07230036551999
$ws.Range("S2").Value = 0.01420895420087014
$ws.Range("T2").Value = 0.01420895420087014

# Row 3
$ws.Range("G3").Value = 4.019613333333333
$ws.Range("H3").Value = 12.05884
$ws.Range("I3").Value = 0.4377217086785624
$ws.Range("J3").Value = 0.4377217086785624
$ws.Range("O3").Value = 0.9174921121349238
$ws.Range("P3").Value = 0.9174921121349238
$ws.Range("Q3").Value = 150.9700176405778
$ws.Range("R3").Value = 1358.7301587652
$ws.Range("S3").Value = 0.401606215022802
$ws.Range("T3").Value = 0.401606215022802

# Row 4
$ws.Range("G4").Value = 4.019613333333333
$ws.Range("H4").Value = 12.05884
$ws.Range("I4").Value = 0.4377217086785624
$ws.Range("J4").Value = 0.4377217086785624
$ws.Range("M4").Value = 2.048706666666666
$ws.Range("N4").Value = 6.14612
$ws.Range("O4").Value = 0.05004672836771996
$ws.Range("P4").Value = 0.05004672836771996
$ws.Range("Q4").Value = 8.23500863342222
$ws.Range("R4").Value = 74.11507770079999
$ws.Range("S4").Value = 0.02190653945489026
$ws.Range("T4").Value = 0.02190653945489026

# Row 5
$ws.Range("I5").Value = 0.02575118419467902
$ws.Range("J5").Value = 0.02575118419467902
$ws.Range("M5").Value = 1.328826
$ws.Range("N5").Value = 3.986478
$ws.Range("O5").Value = 0.03246115949735631
$ws.Range("P5").Value = 0.03246115949735631
$ws.Range("Q5").Value = 0.314232799524
$ws.Range("R5").Value = 2.828095195716
$ws.Range("S5").Value = 0.0008359132973892763
$ws.Range("T5").Value = 0.0008359132973892763

# Row 6
$ws.Range("I6").Value = 0.02575118419467902
$ws.Range("J6").Value = 0.02575118419467902
$ws.Range("O6").Value = 0.9174921121349238
$ws.Range("P6").Value = 0.9174921121349238
$ws.Range("S6").Value = 0.02362650837675152
$ws.Range("T6").Value = 0.02362650837675152

# Row 7
$ws.Range("I7").Value = 0.02575118419467902
$ws.Range("J7").Value = 0.02575118419467902
$ws.Range("M7").Value = 2.048706666666666
$ws.Range("N7").Value = 6.14612
$ws.Range("O7").Value = 0.05004672836771996
$ws.Range("P7").Value = 0.05004672836771996
$ws.Range("Q7").Value = 0.4844658602933333
$ws.Range("R7").Value = 4.36019274264
$ws.Range("S7").Value = 0.001288762520538224
$ws.Range("T7").Value = 0.001288762520538224

# Row 8
$ws.Range("G8").Value = 4.926946666666667
$ws.Range("H8").Value = 14.78084
$ws.Range("I8").Value = 0.5365271071267587
$ws.Range("J8").Value = 0.5365271071267587
$ws.Range("M8").Value = 1.328826
$ws.Range("N8").Value = 3.986478
$ws.Range("O8").Value = 0.03246115949735631
$ws.Range("P8").Value = 0.03246115949735631
$ws.Range("Q8").Value = 6.547054831280001
$ws.Range("R8").Value = 58.92349348152
$ws.Range("S8").Value = 0.01741629199909689
$ws.Range("T8").Value = 0.01741629199909689

# Row 9
$ws.Range("G9").Value = 4.926946666666667
$ws.Range("H9").Value = 14.78084
$ws.Range("I9").Value = 0.5365271071267587
$ws.Range("J9").Value = 0.5365271071267587
$ws.Range("O9").Value = 0.9174921121349238
$ws.Range("P9").Value = 0.9174921121349238
$ws.Range("Q9").Value = 185.0479544916889
$ws.Range("R9").Value = 1665.4315904252
$ws.Range("S9").Value = 0.4922593887353703
$ws.Range("T9").Value = 0.4922593887353703

# Row 10
$ws.Range("G10").Value = 4.926946666666667
$ws.Range("H10").Value = 14.78084
$ws.Range("I10").Value = 0.5365271071267587
$ws.Range("J10").Value = 0.5365271071267587
$ws.Range("M10").Value = 2.048706666666666
$ws.Range("N10").Value = 6.14612
$ws.Range("O10").Value = 0.05004672836771996
$ws.Range("P10").Value = 0.05004672836771996
$ws.Range("Q10").Value = 10.09386848231111
$ws.Range("R10").Value = 90.84481634080001
$ws.Range("S10").Value = 0.02685142639229148
$ws.Range("T10").Value = 0.02685142639229148
